$d = $word.ActiveDocument

# Remove the trailing space after "Welcome to Git"
$d.Content.Find.Execute("Welcome to Git ", $true, $false, $false, $false, $false, $true, 1, $false, "Welcome to Git", 2)

# Add a new paragraph with the commit message text
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("I have done my first commit today")

# Add a new paragraph containing a single space
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter(" ")
